# Updating postprocessing code and files
#
# Applies the numeric corrections from the "Model Fit Tables" document.
#
# NOTE on this runtime: `$d.Tables.Item(N)` handles are not independent
# snapshots. If a table handle captured in a variable is read/written
# *after* a different table index has since been accessed, the stale
# handle silently re-seats onto that other table. To stay safe we always
# resolve `$d.Tables.Item($tableIndex)` freshly, immediately before each
# Cell() access, instead of caching a table reference across calls that
# touch more than one table.

$d = $word.ActiveDocument

function Set-CellValue($tableIndex, $row, $col, $oldValue, $newValue) {
    $current = $d.Tables.Item($tableIndex).Cell($row, $col).Range.Text
    $current = $current.Replace([char]7, "").Replace([char]13, "")
    if ($current -ne $oldValue) {
        Write-Output ("WARNING: Table $tableIndex Cell($row,$col) expected [$oldValue] but found [$current]")
    }
    $d.Tables.Item($tableIndex).Cell($row, $col).Range.Text = $newValue
}

# --- Table 1 : per-period R2 / RMSE / RMSE-SD model fit table ---

# Ice Break row (Period 2 R2/RMSE/RMSE-SD, Period 3 R2/RMSE)
Set-CellValue 1 3 5  "0.79"  "0.78"
Set-CellValue 1 3 6  "3.33"  "3.46"
Set-CellValue 1 3 7  "0.42"  "0.44"
Set-CellValue 1 3 8  "0.70"  "0.88"
Set-CellValue 1 3 9  "5.73"  "5.76"

# Ice Freeze row (Period 3 R2/RMSE/RMSE-SD)
Set-CellValue 1 4 8  "0.45"  "0.44"
Set-CellValue 1 4 9  "5.88"  "5.93"
Set-CellValue 1 4 10 "0.71"  "0.72"

# Temp 1 m row (Period 2 R2/RMSE/RMSE-SD, Period 3 RMSE/RMSE-SD)
Set-CellValue 1 5 5  "0.94"  "0.93"
Set-CellValue 1 5 6  "1.41"  "1.44"
Set-CellValue 1 5 7  "0.25"  "0.26"
Set-CellValue 1 5 9  "2.05"  "2.07"
Set-CellValue 1 5 10 "0.34"  "0.35"

# Temp 4 m row (Period 2 R2/RMSE, Period 3 RMSE)
Set-CellValue 1 6 5  "0.41"  "0.42"
Set-CellValue 1 6 6  "1.83"  "1.82"
Set-CellValue 1 6 9  "1.69"  "1.70"

# Temp 9 m row (Period 3 RMSE/SD)
Set-CellValue 1 7 7  "0.88"  "0.89"

# O2 4 m row (Period 2 RMSE)
Set-CellValue 1 8 6  "4.61"  "5.08"

# TDP row (Period 2 RMSE/RMSE-SD)
Set-CellValue 1 10 6 "3.60"  "3.64"
Set-CellValue 1 10 7 "0.99"  "1.00"

# PP row (Period 2 RMSE, Period 3 R2)
Set-CellValue 1 11 6 "21.26" "21.23"
Set-CellValue 1 11 8 "0.00"  "0.01"

# --- Table 2 : Cumulative PP Modeled/Observed summary ---

# Period 2 row
Set-CellValue 2 3 2  "0.87"  "0.97"
Set-CellValue 2 3 3  "0.93"  "1.01"
Set-CellValue 2 3 4  "0.31"  "0.32"
Set-CellValue 2 3 5  "0.43"  "0.51"
Set-CellValue 2 3 6  "1.38"  "1.52"

# Period 3 row
Set-CellValue 2 4 4  "0.23"  "0.22"
